$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Samples")

# Replace the micro symbol in the two volume headers (uL -> µL)
$ws.Range("F1").Value = "Flow Sample Volume (µL)"
$ws.Range("G1").Value = "Flow PBS Volume (µL)"

# Update the active selection on the Samples sheet to G4
$ws.Activate()
$ws.Range("G4").Select()
